# The sheet held 63 rows of regression data (X1..X4, Y) with no header row.
# The commit adds a header row of column labels (X1, X2, X3, X4, Y) above the
# existing data (pushing every data row down by one, rows 1-63 -> rows 2-64)
# and appends one brand-new trailing data row (new row 64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; all existing data (old rows 1-63) shifts
# down to become rows 2-64, already carrying the correct values/styles - this
# is exactly how the new trailing row 64 (old row 63's data) appears too.
$ws.Rows.Item(1).Insert()

# The new header row inherited the default (no) style from the insert; copy
# formatting from the row below so it matches the rest of the sheet (style
# index 1: center aligned / General number format) for A1:D1, and from
# column A (also style index 1) into E1 specifically, since column E normally
# carries a different number-formatted style (index 2) that the header should
# not use.
$ws.Range("A2:D2").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# Fill in the header labels.
$ws.Range("A1").Value = "X1"
$ws.Range("B1").Value = "X2"
$ws.Range("C1").Value = "X3"
$ws.Range("D1").Value = "X4"
$ws.Range("E1").Value = "Y"

# Restore the saved view state: cursor/selection on H8, scrolled back to the
# top-left (the original file had topLeftCell="A10" / selection K18).
$ws.Range("H8").Select()

Write-Output "done"
